# Actualización del sitio Quarto
# Adds 4 new survey/database entries to the "encuestas" sheet and
# refreshes the active-sheet/selection view state.

$wb = $excel.ActiveWorkbook
$wsEncuestas = $wb.Worksheets.Item("encuestas")

# --- New data rows (31-34) on "encuestas" ---------------------------------

# Row 31: BEEPS (EBRD & World Bank)
$wsEncuestas.Range("A31").Value = "1999-2025"
$wsEncuestas.Range("B31").Value = "EBDR"
$wsEncuestas.Range("C31").Value = "Business Environment and Enterprise Performance Survey (BEEPS) – EBRD & World Bank"
$wsEncuestas.Range("D31").Value = "https://www.ebrd.com/home/what-we-do/office-of-the-chief-economist/beeps/business-environment-enterprise-performance-survey-data.html"
$wsEncuestas.Range("E31").Value = "EXCEL/CSV/STATA"
$wsEncuestas.Range("F31").Value = "Levanta información estandarizada sobre el entorno de negocios, desempeño de las empresas, obstáculos a la inversión, finanzas, innovación, infraestructura y gobernanza. Está diseñada para comparar países y medir la evolución del clima empresarial, siendo una de las principales fuentes internacionales para estudios de productividad e innovación en economías en transición y mercados emergentes"

# Row 32: EDIT (DANE, Colombia)
$wsEncuestas.Range("A32").Value = "2021-2025"
$wsEncuestas.Range("B32").Value = "DANE"
$wsEncuestas.Range("C32").Value = "Encuesta de Desarrollo e Innovación Tecnológica (EDIT) – DANE (Colombia)"
$wsEncuestas.Range("D32").Value = "https://www.dane.gov.co/index.php/estadisticas-por-tema/tecnologia-e-innovacion/encuesta-de-desarrollo-e-innovacion-tecnologica-edit"
$wsEncuestas.Range("E32").Value = "EXCEL/CSV/STATA"
$wsEncuestas.Range("F32").Value = "Recoge datos detallados sobre actividades de innovación, inversión en tecnología, adopción digital, barreras, resultados innovadores y vínculos con universidades y centros tecnológicos. Se ajusta a los lineamientos del Manual de Oslo, permitiendo comparaciones internacionales y análisis de capacidades tecnológicas en el sector productivo colombiano."

# Row 33: Encuesta Nacional de Innovación (PRODUCE, Perú)
$wsEncuestas.Range("A33").Value = "S.F"
$wsEncuestas.Range("B33").Value = "PRODUCE"
$wsEncuestas.Range("C33").Value = "Encuesta Nacional de Innovación – PRODUCE (Perú)"
$wsEncuestas.Range("D33").Value = "https://ogeiee.produce.gob.pe/index.php/en/shortcode/normatividad-metodologia-oee/encuesta-nacional-de-innovacion"
$wsEncuestas.Range("E33").Value = "EXCEL/CSV/STATA"
$wsEncuestas.Range("F33").Value = "Mide las actividades de innovación, capacidades tecnológicas, inversión en I+D, impactos productivos y desafíos que enfrentan las empresas. Su metodología sigue estándares OCDE, proporcionando datos fundamentales para evaluar la competitividad del sector empresarial peruano, así como la adopción tecnológica y el desempeño innovador por sectores y tamaños de firma"

# Row 34: ACTI (Ecuador)
$wsEncuestas.Range("A34").Value = "2012-2014"
$wsEncuestas.Range("B34").Value = "INCEC"
$wsEncuestas.Range("C34").Value = "Encuesta Nacional de Actividades de Ciencia, Tecnología e Innovación (ACTI) – Ecuador"
$wsEncuestas.Range("D34").Value = "https://www.ecuadorencifras.gob.ec/encuesta-nacional-de-actividades-de-ciencia-tecnologia-e-innovacion-acti/"
$wsEncuestas.Range("E34").Value = "EXCEL/CSV"
$wsEncuestas.Range("F34").Value = "Constituye la principal fuente para analizar el ecosistema de ciencia, tecnología e innovación en Ecuador, proporcionando series comparables con estándares internacionales y utilizadas para diagnósticos nacionales y estudios académicos."

# Match the "Siglas" column formatting (column B) used by the rest of the
# table (black font, style index 2) by copying it from the row above.
$wsEncuestas.Range("B30").Copy() | Out-Null
$wsEncuestas.Range("B31:B34").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- View / selection state -------------------------------------------

# "encuestas" becomes the active (tabSelected) sheet, scrolled down toward
# the newly added rows, with the selection left on F38.
$wsEncuestas.Activate()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$wsEncuestas.Range("F38").Select() | Out-Null

Write-Output "Added rows 31-34 to 'encuestas' and updated active sheet/selection"
